$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, copying style from G1 (bold, bordered, centered header style)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill H2:H7 with 0
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
